$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at 38, pushing the existing rows 38-43 down to 42-47.
$ws.Rows("38:41").Insert()

# Constant values shared by every data row in this sheet.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad  = "Cultivar IV Región"
$unidad    = "$/kilo (en caja de 15 kilos)"
$origen    = "Provincia de Limarí"
$kgUnidad  = 1

# New week of data (fecha 44476) for the four quality grades.
$newRows = @(
    @{ Row = 38; Calidad = "Especial";                 Volumen = 360; Min = 2400; Max = 2500; Prom = 2450; PrecioKg = 2450 },
    @{ Row = 39; Calidad = "Extra (doble especial)";    Volumen = 300; Min = 2700; Max = 2800; Prom = 2750; PrecioKg = 2750 },
    @{ Row = 40; Calidad = "Primera";                   Volumen = 500; Min = 2100; Max = 2200; Prom = 2150; PrecioKg = 2150 },
    @{ Row = 41; Calidad = "Segunda";                   Volumen = 400; Min = 1600; Max = 1700; Prom = 1650; PrecioKg = 1650 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = 44476
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.Min
    $ws.Cells.Item($row, 15).Value2 = $r.Max
    $ws.Cells.Item($row, 16).Value2 = $r.Prom
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}
